$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 324, shifting existing rows 324:362 down to 325:363.
$ws.Rows("324:324").Insert()

# Populate the newly inserted row 324 with the new record's data.
$ws.Cells.Item(324, 1).Value2 = 7
$ws.Cells.Item(324, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(324, 3).Value2 = "Ñuble"
$ws.Cells.Item(324, 4).Value2 = 45142
$ws.Cells.Item(324, 5).Value2 = 16
$ws.Cells.Item(324, 6).Value2 = 100112043
$ws.Cells.Item(324, 7).Value2 = "Pepino ensalada"
$ws.Cells.Item(324, 8).Value2 = "Sin especificar"
$ws.Cells.Item(324, 9).Value2 = "Primera"
$ws.Cells.Item(324, 10).Value2 = 60
$ws.Cells.Item(324, 11).Value2 = 12000
$ws.Cells.Item(324, 12).Value2 = 12000
$ws.Cells.Item(324, 13).Value2 = 12000
$ws.Cells.Item(324, 14).Value2 = "`$/caja 60 unidades"
$ws.Cells.Item(324, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(324, 16).Value2 = 200
$ws.Cells.Item(324, 17).Value2 = 60
$ws.Cells.Item(324, 18).Value2 = "Hortaliza"
